# Update block-data metadata: add an "Organ ID" / "Organ Description" split
# by inserting a new column for the organ's descriptive name, and make
# block-data the active sheet/cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("block-data")

# Capture column B's width before inserting so the new column can match it.
$organColWidth = $ws.Columns.Item(2).ColumnWidth()

# Insert a new blank column at C (shifts old Order.."Proteomics" block right).
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(3).ColumnWidth = $organColWidth

# Header row: old "Organ" column becomes "Organ ID"; new C column is
# "Organ Description".
$ws.Range("B1").Value = "Organ ID"
$ws.Range("C1").Value = "Organ Description"

# Fill in the organ description for every tissue-block row, keyed off the
# existing Organ ID (P1 / P2) in column B.
$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    $organId = $ws.Cells.Item($r, 2).Value()
    if ($organId -eq "P1") {
        $ws.Cells.Item($r, 3).Value = "Pancreas 1"
    } else {
        $ws.Cells.Item($r, 3).Value = "Pancreas 2"
    }
}

# Make block-data the active sheet with C12 selected (matches the saved
# workbook state after the edit).
$ws.Activate() | Out-Null
$ws.Range("C12").Select() | Out-Null
